# Update "Pais" sheet: refreshed COVID country stats + a couple of
# leaderboard-order swaps, and a refreshed "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 16:19"

# --- Helper: write a full data row (country name + 7 numeric columns) -
# NOTE: positional parameters only - named parameters do not bind
# correctly for user-defined functions in this host.
function Set-CountryRow($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Straightforward numeric refreshes --------------------------------
Set-CountryRow 4  "Estados Unidos" 4175198 4880  1981505 2046300 0 44  147393
Set-CountryRow 5  "Brasil"         2292286 2335  1570237 637798  0 44  84251
Set-CountryRow 6  "India"          1312551 24421 831059  450671  0 176 30821
Set-CountryRow 21 "Alemania"       205294  152   189400  6704    0 3   9190
Set-CountryRow 58 "Azerbaiyan"     29312   332   21547   7365    0 9   400
Set-CountryRow 62 "Serbia"         22852   409   14047   8287    0 10  518
Set-CountryRow 85 "Noruega"        9088    3     8674    159     0 0   255

# --- Leaderboard-order swap: Tayikistan overtakes Guayana Francesa ----
Set-CountryRow 91 "Tayikistan"       7104 44 5851 1195 0 0 58
Set-CountryRow 92 "Guayana Francesa" 7086 0  5376 1670 0 0 40

# --- Leaderboard-order swap: Namibia jumps ahead of Mozambique/Nueva Zelanda
Set-CountryRow 134 "Namibia"       1618 96 72   1539 0 0 7
Set-CountryRow 135 "Mozambique"    1582 0  528  1043 0 0 11
Set-CountryRow 136 "Nueva Zelanda" 1556 1  1513 21   0 0 22

# --- Leaderboard-order swap: Islas Malvinas ahead of Groenlandia (tie on data)
Set-CountryRow 210 "Islas Malvinas" 13 0 13 0 0 0 0
Set-CountryRow 211 "Groenlandia"    13 0 13 0 0 0 0
